$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "State" tag labels in column A to the new "*Hire" names.
# Column B ("Word") values are left untouched; only the state tags change:
#   No       -> NoHire      (rows 10-14)
#   Maybe    -> MaybeHire   (rows 15-19)
#   Yes      -> YesHire     (rows 20-23)
#   Fuck Yes -> HireNow     (rows 24-35)

$renames = @{
    10 = "NoHire"; 11 = "NoHire"; 12 = "NoHire"; 13 = "NoHire"; 14 = "NoHire";
    15 = "MaybeHire"; 16 = "MaybeHire"; 17 = "MaybeHire"; 18 = "MaybeHire"; 19 = "MaybeHire";
    20 = "YesHire"; 21 = "YesHire"; 22 = "YesHire"; 23 = "YesHire";
    24 = "HireNow"; 25 = "HireNow"; 26 = "HireNow"; 27 = "HireNow"; 28 = "HireNow";
    29 = "HireNow"; 30 = "HireNow"; 31 = "HireNow"; 32 = "HireNow"; 33 = "HireNow";
    34 = "HireNow"; 35 = "HireNow";
}

foreach ($row in $renames.Keys) {
    $ws.Cells.Item($row, 1).Value = $renames[$row]
}

# Reflect the updated selection / scroll position of the sheet view
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 4
